$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1139.5555
$ws.Range("J17").Value = 1169.5
$ws.Range("L17").Value = 3508.5
$ws.Range("N17").Value = -3844.5
$ws.Range("H70").Value = 1509.3889
$ws.Range("I70").Value = 1451.3572
$ws.Range("K70").Value = 4354.071599999999
$ws.Range("M70").Value = -4084.071599999999
$ws.Range("H73").Value = 1509.3889
$ws.Range("I73").Value = 1451.3572
$ws.Range("K73").Value = 4354.071599999999
$ws.Range("M73").Value = -3418.071599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1677.8
$ws.Range("I45").Value = 1813
$ws.Range("K45").Value = 1813
$ws.Range("M45").Value = -1436

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1973
$ws.Range("I86").Value = 2500
$ws.Range("K86").Value = 2500
$ws.Range("M86").Value = -1377
$ws.Range("H89").Value = 1973
$ws.Range("I89").Value = 2500
$ws.Range("K89").Value = 12500
$ws.Range("M89").Value = -6884
$ws.Range("H94").Value = 2841.2
$ws.Range("I94").Value = 1398.6666
$ws.Range("K94").Value = 1398.6666
$ws.Range("M94").Value = -947.6666
$ws.Range("H105").Value = 1062.1111
$ws.Range("I105").Value = 992.5714
$ws.Range("K105").Value = 992.5714
$ws.Range("M105").Value = 754.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4010
$ws.Range("J16").Value = 4010
$ws.Range("L16").Value = 4010
$ws.Range("N16").Value = -4584
$ws.Range("H62").Value = 9166.5
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376
$ws.Range("H65").Value = 9166.5
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880
$ws.Range("H74").Value = 45000
$ws.Range("J74").Value = 45000
$ws.Range("L74").Value = 45000
$ws.Range("N74").Value = -46748
$ws.Range("H77").Value = 45000
$ws.Range("J77").Value = 45000
$ws.Range("L77").Value = 135000
$ws.Range("N77").Value = -143736
$ws.Range("H113").Value = 4010
$ws.Range("J113").Value = 4010
$ws.Range("L113").Value = 4010
$ws.Range("N113").Value = -8350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11000
$ws.Range("I80").Value = 2900
$ws.Range("J80").Value = 21800
$ws.Range("K80").Value = 2900
$ws.Range("L80").Value = 21800
$ws.Range("M80").Value = -1902
$ws.Range("N80").Value = -23796
$ws.Range("H83").Value = 11000
$ws.Range("I83").Value = 2900
$ws.Range("J83").Value = 21800
$ws.Range("K83").Value = 14500
$ws.Range("L83").Value = 109000
$ws.Range("M83").Value = -9508
$ws.Range("N83").Value = -118984
$ws.Range("H97").Value = 299.75
$ws.Range("I97").Value = 299.75
$ws.Range("K97").Value = 299.75
$ws.Range("M97").Value = 196.25
$ws.Range("H113").Value = 2372.5
$ws.Range("I113").Value = 1858.75
$ws.Range("K113").Value = 1858.75
$ws.Range("M113").Value = 311.25
$ws.Range("H122").Value = 1926.3334
$ws.Range("I122").Value = 1926.3334
$ws.Range("K122").Value = 5779.0002
$ws.Range("M122").Value = -3329.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1508
$ws.Range("I7").Value = 1474.75
$ws.Range("K7").Value = 1474.75
$ws.Range("M7").Value = -1362.75
$ws.Range("H16").Value = 8649.5
$ws.Range("I16").Value = 4866.3335
$ws.Range("J16").Value = 19999
$ws.Range("K16").Value = 4866.3335
$ws.Range("L16").Value = 19999
$ws.Range("M16").Value = -4696.3335
$ws.Range("N16").Value = -20339
$ws.Range("H22").Value = 1313.2858
$ws.Range("I22").Value = 3999
$ws.Range("J22").Value = 865.6667
$ws.Range("K22").Value = 3999
$ws.Range("L22").Value = 865.6667
$ws.Range("M22").Value = -3704
$ws.Range("N22").Value = -1455.6667
$ws.Range("H27").Value = 1313.2858
$ws.Range("I27").Value = 3999
$ws.Range("J27").Value = 865.6667
$ws.Range("K27").Value = 3999
$ws.Range("L27").Value = 865.6667
$ws.Range("M27").Value = -3892
$ws.Range("N27").Value = -1079.6667
$ws.Range("H61").Value = 1080
$ws.Range("J61").Value = 870
$ws.Range("L61").Value = 870
$ws.Range("N61").Value = -1274
$ws.Range("H82").Value = 2135.5
$ws.Range("I82").Value = 1894.5
$ws.Range("J82").Value = 2376.5
$ws.Range("K82").Value = 1894.5
$ws.Range("L82").Value = 2376.5
$ws.Range("M82").Value = -1533.5
$ws.Range("N82").Value = -3098.5
$ws.Range("H85").Value = 2135.5
$ws.Range("I85").Value = 1894.5
$ws.Range("J85").Value = 2376.5
$ws.Range("K85").Value = 1894.5
$ws.Range("L85").Value = 2376.5
$ws.Range("M85").Value = -646.5
$ws.Range("N85").Value = -4872.5
$ws.Range("H100").Value = 1166.6666
$ws.Range("I100").Value = 1166.6666
$ws.Range("K100").Value = 1166.6666
$ws.Range("M100").Value = -625.6666
$ws.Range("H113").Value = 1080
$ws.Range("J113").Value = 870
$ws.Range("L113").Value = 870
$ws.Range("N113").Value = -5210
$ws.Range("H126").Value = 1508
$ws.Range("I126").Value = 1474.75
$ws.Range("K126").Value = 4424.25
$ws.Range("M126").Value = -1954.25
$ws.Range("H132").Value = 5329
$ws.Range("I132").Value = 5329
$ws.Range("K132").Value = 15987
$ws.Range("M132").Value = -13457
$ws.Range("H136").Value = 96318.21000000001
$ws.Range("I136").Value = 69864.17999999999
$ws.Range("J136").Value = 193316.33
$ws.Range("K136").Value = 209592.54
$ws.Range("L136").Value = 579948.99
$ws.Range("M136").Value = -207042.54
$ws.Range("N136").Value = -585048.99

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 50000
$ws.Range("J46").Value = 50000
$ws.Range("L46").Value = 50000
$ws.Range("N46").Value = -50462
$ws.Range("H62").Value = 5994
$ws.Range("J62").Value = 5994
$ws.Range("L62").Value = 5994
$ws.Range("N62").Value = -7242
$ws.Range("H65").Value = 5994
$ws.Range("J65").Value = 5994
$ws.Range("L65").Value = 29970
$ws.Range("N65").Value = -36210
$ws.Range("H96").Value = 3880
$ws.Range("I96").Value = 4000
$ws.Range("J96").Value = 3800
$ws.Range("K96").Value = 4000
$ws.Range("L96").Value = 3800
$ws.Range("M96").Value = -2627
$ws.Range("N96").Value = -6546
$ws.Range("H107").Value = 1735.5
$ws.Range("I107").Value = 756.5
$ws.Range("K107").Value = 2269.5
$ws.Range("M107").Value = -349.5
$ws.Range("H126").Value = 2000
$ws.Range("J126").Value = 2000
$ws.Range("L126").Value = 6000
$ws.Range("N126").Value = -10940
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -155070
